$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 3 to make room for the new DEMOG variables
# (NOM, NIM & Natural Increase), pushing the remaining rows down.
$ws.Rows("3:6").Insert()

# Column A (mnemonics) for the 4 new rows
$ws.Range("A3").Value = "NIMTOT"
$ws.Range("A4").Value = "NOMTOT"
$ws.Range("A5").Value = "NATOT"

# Column B (attribute/description) for the first 3 new rows
$ws.Range("B3").Value = "Net Interstate Migration"
$ws.Range("B4").Value = "Net Overseas Migration"
$ws.Range("B5").Value = "Natural Increase"

# Row 6 (POPINC) added after the others
$ws.Range("A6").Value = "POPINC"
$ws.Range("B6").Value = "Total Population Increase"

# Column C (unit) for all 4 new rows matches the Population row above ('000s)
$ws.Range("C3").Value = "('000s)"
$ws.Range("C4").Value = "('000s)"
$ws.Range("C5").Value = "('000s)"
$ws.Range("C6").Value = "('000s)"

# Re-apply the centered number formatting block that spans H27:L32 in the
# final layout (incidental formatting left over on the GVA rows block).
$ws.Range("H27:L30").NumberFormat = "0"
$ws.Range("H27:L30").HorizontalAlignment = -4108

$ws.Range("H31").NumberFormat = "0"
$ws.Range("H31").HorizontalAlignment = -4108
$ws.Range("I31").NumberFormat = "0"
$ws.Range("I31").HorizontalAlignment = -4108
$ws.Range("K31").NumberFormat = "0"
$ws.Range("K31").HorizontalAlignment = -4108
$ws.Range("L31").NumberFormat = "0"
$ws.Range("L31").HorizontalAlignment = -4108

$ws.Range("I32").NumberFormat = "0"
$ws.Range("I32").HorizontalAlignment = -4108
$ws.Range("L32").NumberFormat = "0"
$ws.Range("L32").HorizontalAlignment = -4108

# Move the active selection to D6, matching the saved cursor position.
[void]$ws.Range("D6").Select()
